$d = $word.ActiveDocument

# Use Selection.TypeText on a precisely-bounded Range (paragraph text minus
# the trailing paragraph mark) instead of Find.Execute's ReplaceWith for
# every replacement below. Find.Execute's replacement text silently gets
# "smart quoted" (straight ' -> curly ') by this runtime's autoformatting,
# which corrupts "game's" in the meta description; driving the edit through
# Selection.TypeText writes the literal characters we pass in.

function Replace-ParagraphText($oldText, $newText) {
    $hits = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -eq ($oldText + "`r")) {
            $r = $d.Range($p.Range.Start, $p.Range.End - 1)
            $r.Select()
            $word.Selection.TypeText($newText)
            $hits = $hits + 1
        }
    }
    return $hits
}

# 1. Title heading: remove " - Review" suffix (appears twice: Heading1 and
#    the bold run near the bottom of the document).
Replace-ParagraphText "Play Delicious Candy PopWins for Free – Review" "Play Delicious Candy PopWins for Free"

# 2. "What we like" bullet list items
Replace-ParagraphText "Expanding grid increases ways to win" "Colorful candy symbols create an appealing visual experience"
Replace-ParagraphText "Unique playing grid" "Unique grid expansion mechanic adds excitement to gameplay"
Replace-ParagraphText "Free Spins with multiplier feature" "Free Spins and Multiplier Wheel offer the potential for big wins"
Replace-ParagraphText "Playable on all devices" "Accessible gameplay on all devices with HTML5 technology"

# 3. "What we don't like" bullet list items
Replace-ParagraphText "High volatility" "High volatility may result in less frequent wins"
Replace-ParagraphText "Average RTP rate" "Factor multiplier feature can be randomly activated and may affect gameplay"

# 4. Meta description (italic run)
Replace-ParagraphText "Read our review of Delicious Candy PopWins. Play this unique expanding grid game for free and win big during Free Spins with multiplier feature." "Read a review of Delicious Candy PopWins and play for free. Discover the game's features and gameplay mechanics."
